$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Timestamp (A), new Notified Production MW (B)
$data = @(
    @(2, 46056.01041666666, 0),
    @(3, 46056.02083333334, 0),
    @(4, 46056.03125, 0),
    @(5, 46056.04166666666, 0),
    @(6, 46056.05208333334, 0.391),
    @(7, 46056.0625, 0),
    @(8, 46056.07291666666, 0),
    @(9, 46056.08333333334, 0),
    @(10, 46056.09375, 0),
    @(11, 46056.10416666666, 0),
    @(12, 46056.11458333334, 0),
    @(13, 46056.125, 0),
    @(14, 46056.13541666666, 0.551),
    @(15, 46056.14583333334, 0),
    @(16, 46056.15625, 0),
    @(17, 46056.16666666666, 0),
    @(18, 46056.17708333334, 0),
    @(19, 46056.1875, 0),
    @(20, 46056.19791666666, 0),
    @(21, 46056.20833333334, 0),
    @(22, 46056.21875, 0.729),
    @(23, 46056.22916666666, 0.579),
    @(24, 46056.23958333334, 0.777),
    @(25, 46056.25, 0.911),
    @(26, 46056.26041666666, 3.888),
    @(27, 46056.27083333334, 6.269),
    @(28, 46056.28125, 10.33),
    @(29, 46056.29166666666, 27.987),
    @(30, 46056.30208333334, 94.89),
    @(31, 46056.3125, 153.728),
    @(32, 46056.32291666666, 223.974),
    @(33, 46056.33333333334, 304.475),
    @(34, 46056.34375, 495.873),
    @(35, 46056.35416666666, 590.191),
    @(36, 46056.36458333334, 695.13),
    @(37, 46056.375, 790.258),
    @(38, 46056.38541666666, 991.528),
    @(39, 46056.39583333334, 1074.634),
    @(40, 46056.40625, 1150.89),
    @(41, 46056.41666666666, 1209.175),
    @(42, 46056.42708333334, 1354.075),
    @(43, 46056.4375, 1406.88),
    @(44, 46056.44791666666, 1441.141),
    @(45, 46056.45833333334, 1467.619),
    @(46, 46056.46875, 1491.366),
    @(47, 46056.47916666666, 1499.826),
    @(48, 46056.48958333334, 1495.935),
    @(49, 46056.5, 1486.942),
    @(50, 46056.51041666666, 1443.342),
    @(51, 46056.52083333334, 1406.749),
    @(52, 46056.53125, 1363.39),
    @(53, 46056.54166666666, 1305.416),
    @(54, 46056.55208333334, 1206.344),
    @(55, 46056.5625, 1135.926),
    @(56, 46056.57291666666, 1060.71),
    @(57, 46056.58333333334, 966.8049999999999),
    @(58, 46056.59375, 806.421),
    @(59, 46056.60416666666, 712.05),
    @(60, 46056.61458333334, 613.97),
    @(61, 46056.625, 516.684),
    @(62, 46056.63541666666, 331.335),
    @(63, 46056.64583333334, 258.318),
    @(64, 46056.65625, 188.612),
    @(65, 46056.66666666666, 136.21),
    @(66, 46056.67708333334, 70.952),
    @(67, 46056.6875, 45.859),
    @(68, 46056.69791666666, 36.114),
    @(69, 46056.70833333334, 28.511),
    @(70, 46056.71875, 25.029),
    @(71, 46056.72916666666, 25.065),
    @(72, 46056.73958333334, 25.133),
    @(73, 46056.75, 24.885),
    @(74, 46056.76041666666, 8.151),
    @(75, 46056.77083333334, 7.651),
    @(76, 46056.78125, 4.651),
    @(77, 46056.79166666666, 3.451),
    @(78, 46056.80208333334, 2.491),
    @(79, 46056.8125, 0),
    @(80, 46056.82291666666, 0),
    @(81, 46056.83333333334, 0),
    @(82, 46056.84375, 2.651),
    @(83, 46056.85416666666, 0),
    @(84, 46056.86458333334, 0),
    @(85, 46056.875, 0.651),
    @(86, 46056.88541666666, 2.551),
    @(87, 46056.89583333334, 0),
    @(88, 46056.90625, 0.551),
    @(89, 46056.91666666666, 0),
    @(90, 46056.92708333334, 0),
    @(91, 46056.9375, 0),
    @(92, 46056.94791666666, 0),
    @(93, 46056.95833333334, 0),
    @(94, 46056.96875, 0),
    @(95, 46056.97916666666, 0),
    @(96, 46056.98958333334, 0),
    @(97, 46057, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
